# Update the "想去人数" (interest count) values in column F for the rows
# that changed between the two data refreshes, on both sheets that carry
# this data table: "展览" (sheet1) and "全部类型" (sheet4).

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 568
    4  = 549
    7  = 42
    11 = 4655
    12 = 4454
    13 = 13
    15 = 21
    16 = 155
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
